# Bjorn combat dialog — add new combat lines, tweak two existing lines.
$d = $word.ActiveDocument

# --- 1) "Prepare to meet milheim!" -> all-caps boast line ----------------
$d.Content.Find.Execute(
    "Prepare to meet milheim!", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "PREPARE TO MEET THE BEST GUN IN THE WORLD MY CHODE!", 2) | Out-Null

# --- 2) "My grandma hits harder than that!!!" -> upper-case --------------
$d.Content.Find.Execute(
    "My grandma hits harder than that!!!", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "MY GRANDMA HITS HARDER THAN THAT!!!", 2) | Out-Null

# --- 3) "BURN, BURN, BURN HAHAHAHAHA" -> flamethrower line ---------------
$d.Content.Find.Execute(
    "BURN, BURN, BURN HAHAHAHAHA", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I WOULD HAVE KILLED YOU BY NOW IF I HAD A FLAMETHROWER!!!", 2) | Out-Null

# --- 4) Insert six new bullet lines after "IM GONNA KILL YOU WITH MY CHODE!!!" ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "IM GONNA KILL YOU WITH MY CHODE!!!") {
        $target = $p
        break
    }
}

$ellipsis = [char]0x2026
$rsquo = [char]0x2019

$newLines = @(
    "TIME TO GO BACK UNDER GROUND!",
    "IT TOOK ME YEARS TO GET HERE FROM MY HOMELAND AND IT FEELS LIKE ITS TAKING YOU THE SAME AMOUNT OF TIME TO KILL ME!!!!",
    "I HAVE SEEN SCUM LIKE YOU IN MY HOMELAND! THE ONLY DIFFERENCE IS THAT THEY WERE STONGER!",
    "IS THAT ALL!",
    "COME ON IM GETTING BORED HERE" + $ellipsis + " FINISH THE JOB AND KILL ME!! Oh wait you can" + $rsquo + "t",
    "Zzzzzz" + $ellipsis + " oh huh? Ugh you are still trying to attack?"
)

$cur = $target
foreach ($line in $newLines) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    $cur.Range.Text = $line
}
